# Update the canal_raster Path entry (row 5, column B) to point at the new
# canal raster file, as per the commit:
# "Added Imam's canal raster in data/Strat4/new_canal_raster.tif"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "data/Strat4/new_canal_raster.tif"

# Move the active selection to F6, matching the author's final cursor
# position when they saved the file.
$ws.Range("F6").Select()
